$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = $origStyle
}

Set-TextValue "D2" "67.622.99"
Set-TextValue "E2" "  -0.30%  "
Set-TextValue "D3" "3.317.18"
Set-TextValue "E3" "  +1.00%  "
Set-TextValue "E4" "  +0.21%  "
Set-TextValue "D5" "576.84"
Set-TextValue "E5" "  -0.42%  "
Set-TextValue "D6" "173.86"
Set-TextValue "E6" "  -3.82%  "
Set-TextValue "E7" "  -0.07%  "
Set-TextValue "D8" "0.589"
Set-TextValue "E8" "  +0.68%  "
Set-TextValue "D9" "3.315.04"
Set-TextValue "E9" "  +1.03%  "
Set-TextValue "E10" "  +1.09%  "
Set-TextValue "E11" "  +0.27%  "
Set-TextValue "D12" "45.67"
Set-TextValue "E12" "  -0.05%  "
Set-TextValue "E13" "  -0.52%  "
Set-TextValue "D14" "706.04"
Set-TextValue "E14" "  +2.11%  "
Set-TextValue "D15" "3.860.38"
Set-TextValue "E15" "  +1.19%  "
Set-TextValue "E16" "  +0.04%  "
Set-TextValue "D17" "67.633.49"
Set-TextValue "E17" "  -0.37%  "
Set-TextValue "D19" "3.322.36"
Set-TextValue "E19" "  +0.98%  "
Set-TextValue "E20" "  -0.35%  "
Set-TextValue "D22" "0.888"
Set-TextValue "E22" "  -0.57%  "
Set-TextValue "D23" "5.34"
Set-TextValue "D24" "16.85"
Set-TextValue "E24" "  -2.62%  "
Set-TextValue "D25" "98.09"
Set-TextValue "E25" "  +0.29%  "
Set-TextValue "E26" "  -3.07%  "
Set-TextValue "E27" "  -2.69%  "
Set-TextValue "D28" "9.32"
Set-TextValue "E28" "  -0.83%  "
Set-TextValue "D29" "33.29"
Set-TextValue "E29" "  +1.79%  "
Set-TextValue "D30" "8.46"
Set-TextValue "E30" "  +0.22%  "
Set-TextValue "D31" "7.10"
Set-TextValue "E31" "  +5.45%  "
Set-TextValue "D32" "566.95"
Set-TextValue "E32" "  -3.51%  "
Set-TextValue "D33" "10.93"
Set-TextValue "E33" "  +0.67%  "
Set-TextValue "E34" "  +0.68%  "
Set-TextValue "D35" "57.68"
Set-TextValue "E35" "  +4.16%  "
Set-TextValue "D36" "0.998"
Set-TextValue "E36" "  +0.15%  "
Set-TextValue "D37" "3.704.25"
Set-TextValue "E37" "  -4.80%  "
Set-TextValue "E38" "  -2.14%  "
Set-TextValue "D39" "34.27"
Set-TextValue "E39" "  +5.98%  "
Set-TextValue "E40" "  -0.33%  "
Set-TextValue "E41" "  -2.93%  "
Set-TextValue "E42" "  -1.19%  "
Set-TextValue "D43" "3.29"
Set-TextValue "E43" "  -3.19%  "
Set-TextValue "E44" "  +0.17%  "
Set-TextValue "E45" "  -1.64%  "
Set-TextValue "B46" "VeChain"
Set-TextValue "C46" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D46" "0.0404"
Set-TextValue "E46" "  -2.08%  "
Set-TextValue "B47" "ThetaToken"
Set-TextValue "C47" "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue "D47" "2.68"
Set-TextValue "E47" "  +6.80%  "
Set-TextValue "E48" "  -0.13%  "
Set-TextValue "E49" "  -0.33%  "
Set-TextValue "E50" "  -5.33%  "
Set-TextValue "D51" "128.42"
Set-TextValue "E51" "  -1.48%  "
